# Generate Report for handoff
# Mark the defb06e1-... file as "Ready for handoff" across sheets, and
# update the zh-cn / de-de "Latest Handoff Datetime" for that file.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to defb06e1-...md
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3 corresponds to defb06e1-...md
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-01-26 09:39:21"

# de-de sheet: row 3 corresponds to defb06e1-...md
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-01-26 09:39:33"
